# Edit: "The Dragon Boat Festival"
# 1) Replace the paragraphs spanning "2023年3月4日" and the following
#    "今天周末没有课..." paragraph with their updated formatting/text, plus
#    two brand-new trailing paragraphs (a new date line and a new diary
#    entry), preserving the trailing _GoBack bookmark at the very end.
$d = $word.ActiveDocument

$p3 = $d.Paragraphs(3)
$p4 = $d.Paragraphs(4)
$rng = $d.Range($p3.Range.Start, $p4.Range.End)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2023年3月4日</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>今天周末没有课 可以睡到自然醒了 大晴天下午可以晒个太阳 耶！！</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>2023年3月5日</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:rFonts w:hint="default"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="eastAsia"/><w:lang w:val="en-US" w:eastAsia="zh-CN"/></w:rPr><w:t>今天还是周末 继续睡到自然醒 但是明天又要上课了 又是早八！</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@

$rng.InsertXML($xml)

# 2) Styles: mark the "Normal Table" style as a QuickStyle (adds <w:qFormat/>
#    to its style definition, i.e. w:styleId="2").
$tableStyle = $d.Styles("Normal Table")
$tableStyle.QuickStyle = $true
